{"js": "// The commit removes the leading \"Professor \" from the author credit line\n// that appears in the document's (default/\"Primary\") footer, changing:\n//   \"Professor Rui Brito / Lu\u00eds Encerrabodes \"\n// to:\n//   \"Rui Brito / Lu\u00eds Encerrabodes \"\n\nconst sections = context.document.sections;\nsections.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < sections.items.length; i++) {\n  const section = sections.items[i];\n  const footer = section.getFooter(\"Primary\");\n  const hits = footer.search(\"Professor Rui Brito / Lu\u00eds Encerrabodes\", { matchCase: true });\n  hits.load(\"items\");\n  await context.sync();\n\n  for (let j = 0; j < hits.items.length; j++) {\n    hits.items[j].insertText(\"Rui Brito / Lu\u00eds Encerrabodes\", \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# The commit removes the leading \"Professor \" from the author credit line\n# in the document footer, changing:\n#   \"Professor Rui Brito / Lu\u00eds Encerrabodes \"\n# to:\n#   \"Rui Brito / Lu\u00eds Encerrabodes \"\n\n$d = $word.ActiveDocument\n\nforeach ($sec in $d.Sections) {\n    # wdHeaderFooterPrimary=1, wdHeaderFooterFirstPage=2, wdHeaderFooterEvenPages=3\n    for ($i = 1; $i -le 3; $i++) {\n        $footer = $sec.Footers.Item($i)\n        $find = $footer.Range.Find\n        $find.ClearFormatting()\n        $find.Text = \"Professor Rui Brito / Lu\u00eds Encerrabodes\"\n        $find.Replacement.Text = \"Rui Brito / Lu\u00eds Encerrabodes\"\n        $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n    }\n}\n"}
